$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.420.93'
$ws.Range('E2').Value = '  -1.94%  '
$ws.Range('D3').Value = '3.379.94'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.74'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.01'
$ws.Range('E6').Value = '  -6.23%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.380.12'
$ws.Range('E8').Value = '  -1.94%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.472'
$ws.Range('E11').Value = '  -3.16%  '
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '3.957.05'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('E14').Value = '  +0.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.95'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '3.382.53'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('E17').Value = '  -3.27%  '
$ws.Range('D18').Value = '60.565.73'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.18'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.79'
$ws.Range('E20').Value = '  -4.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.99'
$ws.Range('E21').Value = '  -5.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '386.12'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.92'
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000116'
$ws.Range('E26').Value = '  -7.81%  '
$ws.Range('D27').Value = '3.519.53'
$ws.Range('E27').Value = '  -1.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.179'
$ws.Range('E28').Value = '  -1.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.37'
$ws.Range('E30').Value = '  -4.76%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.14'
$ws.Range('E31').Value = '  -2.20%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.90'
$ws.Range('E32').Value = '  -4.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.41'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.40'
$ws.Range('E35').Value = '  -2.79%  '
$ws.Range('D36').Value = '3.410.92'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '168.00'
$ws.Range('E37').Value = '  +1.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.86'
$ws.Range('E38').Value = '  -2.49%  '
$ws.Range('E39').Value = '  -4.88%  '
$ws.Range('E40').Value = '  -4.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0769'
$ws.Range('E41').Value = '  -2.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '27.08'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  -1.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.39'
$ws.Range('E45').Value = '  -2.49%  '
$ws.Range('E46').Value = '  -2.22%  '
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('D48').Value = '2.518.18'
$ws.Range('E48').Value = '  -3.29%  '
$ws.Range('E49').Value = '  -4.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.05'
$ws.Range('E50').Value = '  -0.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.73'
